$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: column K was "Tendencia" (with "scrapy_datetime" in L). The
# "Tendencia" column is dropped and "scrapy_datetime" now lives in K;
# column L is removed entirely.
$ws.Range("K1").Value = "scrapy_datetime"
$ws.Columns.Item(12).Delete()

# Data rows 2-13: these finished scraping a category and now carry real
# counts/prices/trend links instead of placeholder zeros / "NA".
$data = @{
    2  = @("247 resultado", "105 resultado", 129,     159, 262.67,  155, "https://trends.google.com.br/trends/explore?geo=BR&q=meia beach tennis")
    3  = @("1.866 resultado", "15 resultado", 29,      29, 1205,    695, "https://trends.google.com.br/trends/explore?geo=BR&q=bandana tubular")
    4  = @("897 resultado", "17 resultado", 122,      155, 803,     228, "https://trends.google.com.br/trends/explore?geo=BR&q=roupa camuflada")
    5  = @("893 resultado", "9 resultado", 1039.67,   737, 1196,   1639, "https://trends.google.com.br/trends/explore?geo=BR&q=bicicleta triciclo")
    6  = @("8.220 resultado", "148 resultado", 122,    139, 1215,   514, "https://trends.google.com.br/trends/explore?geo=BR&q=camelbak")
    7  = @("2.382 resultado", "36 resultado", 79,       58, 1627.33,1110, "https://trends.google.com.br/trends/explore?geo=BR&q=bermuda termica masculina")
    8  = @("679 resultado", "NaoTem", 526,            579, 469,     150, "https://trends.google.com.br/trends/explore?geo=BR&q=arpao pesca")
    9  = @("42.138 resultado", "401 resultado", 275.67,300, 199,     34, "https://trends.google.com.br/trends/explore?geo=BR&q=kit sobrevivencia completo")
    10 = @("2.830 resultado", "32 resultado", 109,      64, 236,     197, "https://trends.google.com.br/trends/explore?geo=BR&q=garrafa personalizada")
    11 = @("1.751 resultado", "1 resultado", 80,        42, 149.67,  172, "https://trends.google.com.br/trends/explore?geo=BR&q=quechua")
    12 = @("7.236 resultado", "132 resultado", 29.33,   28, 1458,   1920, "https://trends.google.com.br/trends/explore?geo=BR&q=balaclava")
    13 = @("17.508 resultado", "259 resultado", 39.67,  22, 32307, 13580, "https://trends.google.com.br/trends/explore?geo=BR&q=colchonete")
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]   # D: Qnt-Normal
    $ws.Cells.Item($row, 5).Value = $vals[1]   # E: Qnt-FULL
    $ws.Cells.Item($row, 6).Value = $vals[2]   # F: Media-Preco
    $ws.Cells.Item($row, 7).Value = $vals[3]   # G: Mediana-Preco
    $ws.Cells.Item($row, 8).Value = $vals[4]   # H: Media-Vendas
    $ws.Cells.Item($row, 9).Value = $vals[5]   # I: Mediana-Vendas
    $ws.Cells.Item($row, 10).Value = $vals[6]  # J: GoogleTrends
}

# Every row's scrapy_datetime (now column K) is refreshed to the new run.
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 11).Value = "2022-05-21 08:45:30"
}
